$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValue = "UAT2\Antosova\2025\Converged Connectivity\NonRelease\CCCTR-1214 - BB Security balicek"

for ($row = 2; $row -le 141; $row++) {
    $ws.Cells.Item($row, 3).Value = $newValue
}
